$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case connector words in municipality/state names ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B9').Value = 'San Francisco De Los Romo'
$ws.Range('B15').Value = 'Playas De Rosarito'
$ws.Range('B34').Value = 'Amatenango De La Frontera'
$ws.Range('B39').Value = 'Benemérito De Las Américas'
$ws.Range('B46').Value = 'Chiapa De Corzo'
$ws.Range('B51').Value = 'Comitán De Domínguez'
$ws.Range('B73').Value = 'Marqués De Comillas'
$ws.Range('B80').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B89').Value = 'Salto De Agua'
$ws.Range('B90').Value = 'San Cristóbal De Las Casas'
$ws.Range('B132').Value = 'Coyame Del Sotol'
$ws.Range('B143').Value = 'Guadalupe Y Calvo'
$ws.Range('B147').Value = 'Hidalgo Del Parral'
$ws.Range('B169').Value = 'San Francisco De Borja'
$ws.Range('B170').Value = 'San Francisco De Conchos'
$ws.Range('B171').Value = 'San Francisco Del Oro'
$ws.Range('B178').Value = 'Valle De Zaragoza'
$ws.Range('A180').Value = 'Ciudad De México'
$ws.Range('B183').Value = 'Cuajimalpa De Morelos'
$ws.Range('A198').Value = 'Coahuila De Zaragoza'
$ws.Range('B207').Value = 'San Juan De Sabinas'
$ws.Range('B219').Value = 'Villa De Álvarez'
$ws.Range('B223').Value = 'Coneto De Comonfort'
$ws.Range('B238').Value = 'Nombre De Dios'
$ws.Range('B245').Value = 'Pánuco De Coronado'
$ws.Range('B249').Value = 'San Juan Del Río'
$ws.Range('B250').Value = 'San Luis Del Cordero'
$ws.Range('B251').Value = 'San Pedro Del Gallo'
$ws.Range('A261').Value = 'Estado De México'
$ws.Range('B261').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B264').Value = 'Almoloya De Alquisiras'
$ws.Range('B265').Value = 'Almoloya De Juárez'
$ws.Range('B270').Value = 'Atizapán De Zaragoza'
$ws.Range('B276').Value = 'Chapa De Mota'
$ws.Range('B280').Value = 'Coacalco De Berriozábal'
$ws.Range('B286').Value = 'Ecatepec De Morelos'
$ws.Range('B291').Value = 'Ixtapan De La Sal'
$ws.Range('B292').Value = 'Ixtapan Del Oro'
$ws.Range('B300').Value = 'Naucalpan De Juárez'
$ws.Range('B310').Value = 'San Felipe Del Progreso'
$ws.Range('B311').Value = 'San José Del Rincón'
$ws.Range('B312').Value = 'San Martín De Las Pirámides'
$ws.Range('B314').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B322').Value = 'Tenango Del Valle'
$ws.Range('B331').Value = 'Tlalnepantla De Baz'
$ws.Range('B337').Value = 'Valle De Bravo'
$ws.Range('B338').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B341').Value = 'Villa De Allende'
$ws.Range('B342').Value = 'Villa Del Carbón'
$ws.Range('B351').Value = 'Apaseo El Alto'
$ws.Range('B352').Value = 'Apaseo El Grande'
$ws.Range('B361').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B365').Value = 'Jaral Del Progreso'
$ws.Range('B372').Value = 'Purísima Del Rincón'
$ws.Range('B377').Value = 'San Diego De La Unión'
$ws.Range('B379').Value = 'San Francisco Del Rincón'
$ws.Range('B381').Value = 'San Luis De La Paz'
$ws.Range('B382').Value = 'San Miguel De Allende'
$ws.Range('B383').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B384').Value = 'Silao De La Victoria'
$ws.Range('B388').Value = 'Valle De Santiago'
$ws.Range('B394').Value = 'Acapulco De Juárez'
$ws.Range('B397').Value = 'Alcozauca De Guerrero'
$ws.Range('B401').Value = 'Atenango Del Río'
$ws.Range('B402').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B404').Value = 'Atoyac De Álvarez'
$ws.Range('B405').Value = 'Ayutla De Los Libres'
$ws.Range('B408').Value = 'Buenavista De Cuéllar'
$ws.Range('B409').Value = 'Chilapa De Álvarez'
$ws.Range('B410').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B411').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B415').Value = 'Coyuca De Benítez'
$ws.Range('B416').Value = 'Coyuca De Catalán'
$ws.Range('B420').Value = 'Cuetzala Del Progreso'
$ws.Range('B421').Value = 'Cutzamala De Pinzón'
$ws.Range('B427').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B428').Value = 'Iguala De La Independencia'
$ws.Range('B430').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B432').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B436').Value = 'Mártir De Cuilapan'
$ws.Range('B447').Value = 'Taxco De Alarcón'
$ws.Range('B450').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B452').Value = 'Tixtla De Guerrero'
$ws.Range('B455').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B456').Value = 'Tlapa De Comonfort'
$ws.Range('B458').Value = 'Técpan De Galeana'
$ws.Range('B462').Value = 'Zihuatanejo De Azueta'
$ws.Range('B469').Value = 'Agua Blanca De Iturbide'
$ws.Range('B474').Value = 'Atotonilco De Tula'
$ws.Range('B475').Value = 'Atotonilco El Grande'
$ws.Range('B479').Value = 'Cuautepec De Hinojosa'
$ws.Range('B483').Value = 'Huasca De Ocampo'
$ws.Range('B486').Value = 'Huejutla De Reyes'
$ws.Range('B493').Value = 'Mineral De La Reforma'
$ws.Range('B494').Value = 'Mineral Del Chico'
$ws.Range('B495').Value = 'Mineral Del Monte'
$ws.Range('B496').Value = 'Mixquiahuala De Juárez'
$ws.Range('B498').Value = 'Omitlán De Juárez'
$ws.Range('B499').Value = 'Pachuca De Soto'
$ws.Range('B502').Value = 'Progreso De Obregón'
$ws.Range('B506').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B510').Value = 'Tepehuacán De Guerrero'
$ws.Range('B511').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B513').Value = 'Tezontepec De Aldama'
$ws.Range('B519').Value = 'Tula De Allende'
$ws.Range('B520').Value = 'Tulancingo De Bravo'
$ws.Range('B522').Value = 'Zacualtipán De Ángeles'
$ws.Range('B525').Value = 'Acatlán De Juárez'
$ws.Range('B526').Value = 'Ahualulco De Mercado'
$ws.Range('B531').Value = 'Atemajac De Brizuela'
$ws.Range('B534').Value = 'Atotonilco El Alto'
$ws.Range('B536').Value = 'Autlán De Navarro'
$ws.Range('B540').Value = 'Cañadas De Obregón'
$ws.Range('B546').Value = 'Concepción De Buenos Aires'
$ws.Range('B553').Value = 'Encarnación De Díaz'
$ws.Range('B558').Value = 'Huejuquilla El Alto'
$ws.Range('B560').Value = 'Ixtlahuacán Del Río'
$ws.Range('B562').Value = 'Jilotlán De Los Dolores'
$ws.Range('B568').Value = 'Lagos De Moreno'
$ws.Range('B573').Value = 'Ojuelos De Jalisco'
$ws.Range('B578').Value = 'San Cristóbal De La Barranca'
$ws.Range('B579').Value = 'San Diego De Alejandría'
$ws.Range('B581').Value = 'San Juan De Los Lagos'
$ws.Range('B582').Value = 'San Juanito De Escobedo'
$ws.Range('B585').Value = 'San Martín De Bolaños'
$ws.Range('B586').Value = 'San Miguel El Alto'
$ws.Range('B588').Value = 'San Sebastián Del Oeste'
$ws.Range('B589').Value = 'Santa María De Los Ángeles'
$ws.Range('B590').Value = 'Santa María Del Oro'
$ws.Range('B593').Value = 'Talpa De Allende'
$ws.Range('B594').Value = 'Tamazula De Gordiano'
$ws.Range('B599').Value = 'Teocuitatlán De Corona'
$ws.Range('B600').Value = 'Tepatitlán De Morelos'
$ws.Range('B602').Value = 'Tizapán El Alto'
$ws.Range('B603').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B613').Value = 'Unión De San Antonio'
$ws.Range('B614').Value = 'Unión De Tula'
$ws.Range('B619').Value = 'Yahualica De González Gallo'
$ws.Range('B620').Value = 'Zacoalco De Torres'
$ws.Range('B623').Value = 'Zapotitlán De Vadillo'
$ws.Range('B625').Value = 'Zapotlán Del Rey'
$ws.Range('B626').Value = 'Zapotlán El Grande'
$ws.Range('A628').Value = 'Michoacán De Ocampo'
$ws.Range('B647').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B649').Value = 'Cojumatlán De Régules'
$ws.Range('B710').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B744').Value = 'Puente De Ixtla'
$ws.Range('B750').Value = 'Tetela Del Volcán'
$ws.Range('B751').Value = 'Tlaltizapán De Zapata'
$ws.Range('B762').Value = 'Bahía De Banderas'
$ws.Range('B767').Value = 'Ixtlán Del Río'
$ws.Range('B774').Value = 'Santa María Del Oro'
$ws.Range('B789').Value = 'San Nicolás De Los Garza'
$ws.Range('B792').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B799').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B803').Value = 'Cuilápam De Guerrero'
$ws.Range('B804').Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range('B806').Value = 'El Barrio De La Soledad'
$ws.Range('B808').Value = 'Guevea De Humboldt'
$ws.Range('B809').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B810').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B811').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B812').Value = 'Huajuapan De León'
$ws.Range('B813').Value = 'Ixtlán De Juárez'
$ws.Range('B817').Value = 'Mazatlán Villa De Flores'
$ws.Range('B818').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B820').Value = 'Oaxaca De Juárez'
$ws.Range('B821').Value = 'Ocotlán De Morelos'
$ws.Range('B822').Value = 'Putla Villa De Guerrero'
$ws.Range('B829').Value = 'San Antonino El Alto'
$ws.Range('B874').Value = 'San Miguel Del Puerto'
$ws.Range('B876').Value = 'San Pablo Villa De Mitla'
$ws.Range('B884').Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range('B885').Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range('B896').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B897').Value = 'Santa Inés Del Monte'
$ws.Range('B900').Value = 'Santa Lucía Del Camino'
$ws.Range('B907').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B934').Value = 'Santo Domingo De Morelos'
$ws.Range('B938').Value = 'Tataltepec De Valdés'
$ws.Range('B939').Value = 'Teotitlán De Flores Magón'
$ws.Range('B941').Value = 'Tlacolula De Matamoros'
$ws.Range('B942').Value = 'Totontepec Villa De Morelos'
$ws.Range('B943').Value = 'Villa Sola De Vega'
$ws.Range('B944').Value = 'Villa Tejúpam De La Unión'
$ws.Range('B945').Value = 'Villa De Etla'
$ws.Range('B946').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B947').Value = 'Villa De Zaachila'
$ws.Range('B948').Value = 'Zimatlán De Álvarez'
$ws.Range('B960').Value = 'Ayotoxco De Guerrero'
$ws.Range('B962').Value = 'Chalchicomula De Sesma'
$ws.Range('B986').Value = 'Huitzilan De Serdán'
$ws.Range('B988').Value = 'Izúcar De Matamoros'
$ws.Range('B996').Value = 'Los Reyes De Juárez'
$ws.Range('B1001').Value = 'Palmar De Bravo'
$ws.Range('B1014').Value = 'San Salvador El Seco'
$ws.Range('B1015').Value = 'San Salvador El Verde'
$ws.Range('B1025').Value = 'Tepanco De López'
$ws.Range('B1026').Value = 'Tepango De Rodríguez'
$ws.Range('B1030').Value = 'Tepexi De Rodríguez'
$ws.Range('B1031').Value = 'Tetela De Ocampo'
$ws.Range('B1035').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B1057').Value = 'Amealco De Bonfil'
$ws.Range('B1058').Value = 'Cadereyta De Montes'
$ws.Range('B1062').Value = 'Jalpan De Serra'
$ws.Range('B1063').Value = 'Landa De Matamoros'
$ws.Range('B1064').Value = 'Pinal De Amoles'
$ws.Range('B1066').Value = 'San Juan Del Río'
$ws.Range('B1073').Value = 'Armadillo De Los Infante'
$ws.Range('B1085').Value = 'Mexquitic De Carmona'
$ws.Range('B1089').Value = 'San Ciro De Acosta'
$ws.Range('B1094').Value = 'Santa María Del Río'
$ws.Range('B1096').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B1100').Value = 'Villa De Arriaga'
$ws.Range('B1101').Value = 'Villa De Ramos'
$ws.Range('B1165').Value = 'Nacozari De García'
$ws.Range('B1178').Value = 'San Pedro De La Cueva'
$ws.Range('B1199').Value = 'Jalpa De Méndez'
$ws.Range('B1224').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B1225').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B1226').Value = 'San Pablo Del Monte'
$ws.Range('B1230').Value = 'Tetla De La Solidaridad'
$ws.Range('A1237').Value = 'Veracruz De Ignacio De La Llave'
$ws.Range('B1242').Value = 'Amatlán De Los Reyes'
$ws.Range('B1249').Value = 'Boca Del Río'
$ws.Range('B1251').Value = 'Camarón De Tejeda'
$ws.Range('B1254').Value = 'Castillo De Teayo'
$ws.Range('B1256').Value = 'Cazones De Herrera'
$ws.Range('B1268').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1282').Value = 'Hueyapan De Ocampo'
$ws.Range('B1283').Value = 'Ignacio De La Llave'
$ws.Range('B1286').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B1288').Value = 'Ixhuatlán Del Café'
$ws.Range('B1289').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B1298').Value = 'Lerdo De Tejada'
$ws.Range('B1300').Value = 'Martínez De La Torre'
$ws.Range('B1311').Value = 'Paso De Ovejas'
$ws.Range('B1314').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1323').Value = 'Sayula De Alemán'
$ws.Range('B1326').Value = 'Soledad De Doblado'
$ws.Range('B1371').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1375').Value = 'El Plateado De Joaquín Amaro'
$ws.Range('B1386').Value = 'Jiménez Del Teul'
$ws.Range('B1395').Value = 'Moyahua De Estrada'
$ws.Range('B1396').Value = 'Nochistlán De Mejía'
$ws.Range('B1406').Value = 'Teúl De González Ortega'
$ws.Range('B1407').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1412').Value = 'Villa De Cos'
$ws.Range('A1416').Value = 'Total'

# --- Delete footer metadata rows 1418-1422 ---
$ws.Range("A1418:D1422").EntireRow.Delete()

